$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target OOXML column width is 15.42578125 characters. This runtime quantizes
# ColumnWidth to 1/6-character increments when persisting to XML, so the
# nearest reachable stored width is 15.5 (vs. 15.333333333333334 on the other
# side) -- use an input value safely inside that rounding bucket.
$ws.Columns.Item(1).ColumnWidth = 14.65

$values = @(
  0.09282644053402578,
  -0.009999999805774706,
  -0.008999999807953074,
  0.28399664880383213,
  -0.005999999813595558,
  -0.005999999807091427,
  -0.019999999773164134,
  -0.019999999770982768,
  -0.005999999801506561,
  -0.005999999799250588,
  -0.004499999802774823,
  -0.0059999997986595055,
  -0.005999999796507005,
  -0.011999999781711068,
  -0.005999999795576194,
  -0.005999999795331057,
  0.006952780926051894,
  -0.008999999787940638,
  -0.056757752870999933,
  -0.008999999806016845,
  -0.008999999805774372,
  -0.008999999805601178,
  -0.00899999980633126,
  -0.041999999725137016,
  -0.042841264841605664,
  -0.005999999806586942,
  -0.005999999806067358,
  -0.0059999998034188096,
  0.01385292867504262,
  -0.019999999767330134,
  -0.014999999777305817,
  -0.020999999762799426,
  -0.006626268078941955
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $values[$i]
}
